$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "2358"
$t.Cell(5,1).Range.Text  = "0.00003"
$t.Cell(6,1).Range.Text  = "0.00078"
$t.Cell(7,1).Range.Text  = "0.00015"
$t.Cell(9,1).Range.Text  = "0.00026"
$t.Cell(10,1).Range.Text = "0.00030"
$t.Cell(11,1).Range.Text = "0.00036"
$t.Cell(12,1).Range.Text = "0.42080"

$t.Cell(44,1).Range.Text = "99.94"
$t.Cell(45,1).Range.Text = "0.42"
$t.Cell(46,1).Range.Text = "711"
